$wb = $excel.ActiveWorkbook

# Sheet: VENTAS POR GRUPO
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M14").Value = 5102.69

# Sheet: VENTA MENSUAL
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F14").Value = 9712.51
$wsMensual.Range("F23").Value = 66667.3

# Sheet: CUMPLIMIENTO MENSUAL
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D12").Value = 47257.91
$wsCumpl.Range("E12").Value = -2839.910000000003
$wsCumpl.Range("F12").Value = 1.063936016930074
$wsCumpl.Range("D14").Value = 66667.3
$wsCumpl.Range("E14").Value = -11267.82898829906
$wsCumpl.Range("F14").Value = 1.203392357048304
